$d = $word.ActiveDocument

$d.Content.Find.Execute("2025-10-15 Wednesday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-10-16 Thursday", 2) | Out-Null
$d.Content.Find.Execute("21+0=", $true, $false, $false, $false, $false, $true, 1, $false, "58-28=", 2) | Out-Null
$d.Content.Find.Execute("10+26=", $true, $false, $false, $false, $false, $true, 1, $false, "68+14=", 2) | Out-Null
$d.Content.Find.Execute("88-65=", $true, $false, $false, $false, $false, $true, 1, $false, "71-36=", 2) | Out-Null
$d.Content.Find.Execute("99-3=", $true, $false, $false, $false, $false, $true, 1, $false, "86-71=", 2) | Out-Null
$d.Content.Find.Execute("68-56=", $true, $false, $false, $false, $false, $true, 1, $false, "98-66=", 2) | Out-Null
$d.Content.Find.Execute("81-44=", $true, $false, $false, $false, $false, $true, 1, $false, "33+40=", 2) | Out-Null
$d.Content.Find.Execute("5+57=", $true, $false, $false, $false, $false, $true, 1, $false, "92-56=", 2) | Out-Null
$d.Content.Find.Execute("14+24=", $true, $false, $false, $false, $false, $true, 1, $false, "62-62=", 2) | Out-Null
$d.Content.Find.Execute("58-1=", $true, $false, $false, $false, $false, $true, 1, $false, "23+62=", 2) | Out-Null
$d.Content.Find.Execute("45+25=", $true, $false, $false, $false, $false, $true, 1, $false, "1+55=", 2) | Out-Null
$d.Content.Find.Execute("92-52=", $true, $false, $false, $false, $false, $true, 1, $false, "13+35=", 2) | Out-Null
$d.Content.Find.Execute("85-7=", $true, $false, $false, $false, $false, $true, 1, $false, "88-30=", 2) | Out-Null
$d.Content.Find.Execute("47+24=", $true, $false, $false, $false, $false, $true, 1, $false, "33+40=", 2) | Out-Null
$d.Content.Find.Execute("86-60=", $true, $false, $false, $false, $false, $true, 1, $false, "60-17=", 2) | Out-Null
$d.Content.Find.Execute("35+62=", $true, $false, $false, $false, $false, $true, 1, $false, "72-71=", 2) | Out-Null
$d.Content.Find.Execute("15+6=", $true, $false, $false, $false, $false, $true, 1, $false, "90+6=", 2) | Out-Null
$d.Content.Find.Execute("13+52=", $true, $false, $false, $false, $false, $true, 1, $false, "67-6=", 2) | Out-Null
$d.Content.Find.Execute("23-8=", $true, $false, $false, $false, $false, $true, 1, $false, "71-69=", 2) | Out-Null
$d.Content.Find.Execute("47+25=", $true, $false, $false, $false, $false, $true, 1, $false, "41-37=", 2) | Out-Null
$d.Content.Find.Execute("9+88=", $true, $false, $false, $false, $false, $true, 1, $false, "69+15=", 2) | Out-Null
$d.Content.Find.Execute("69-29=", $true, $false, $false, $false, $false, $true, 1, $false, "91+0=", 2) | Out-Null
$d.Content.Find.Execute("80+3=", $true, $false, $false, $false, $false, $true, 1, $false, "17-4=", 2) | Out-Null
$d.Content.Find.Execute("65-30=", $true, $false, $false, $false, $false, $true, 1, $false, "79-59=", 2) | Out-Null
$d.Content.Find.Execute("6+31=", $true, $false, $false, $false, $false, $true, 1, $false, "48+42=", 2) | Out-Null
$d.Content.Find.Execute("8+20=", $true, $false, $false, $false, $false, $true, 1, $false, "54+2=", 2) | Out-Null
$d.Content.Find.Execute("68-27=", $true, $false, $false, $false, $false, $true, 1, $false, "65-10=", 2) | Out-Null
$d.Content.Find.Execute("25+52=", $true, $false, $false, $false, $false, $true, 1, $false, "14+7=", 2) | Out-Null
$d.Content.Find.Execute("60-38=", $true, $false, $false, $false, $false, $true, 1, $false, "51+19=", 2) | Out-Null
$d.Content.Find.Execute("58+18=", $true, $false, $false, $false, $false, $true, 1, $false, "32-6=", 2) | Out-Null
$d.Content.Find.Execute("76-51=", $true, $false, $false, $false, $false, $true, 1, $false, "43-32=", 2) | Out-Null
$d.Content.Find.Execute("48-35=", $true, $false, $false, $false, $false, $true, 1, $false, "27-12=", 2) | Out-Null
$d.Content.Find.Execute("71-45=", $true, $false, $false, $false, $false, $true, 1, $false, "2+61=", 2) | Out-Null
$d.Content.Find.Execute("53-22=", $true, $false, $false, $false, $false, $true, 1, $false, "39-31=", 2) | Out-Null
$d.Content.Find.Execute("69-57=", $true, $false, $false, $false, $false, $true, 1, $false, "52+41=", 2) | Out-Null
$d.Content.Find.Execute("2+86=", $true, $false, $false, $false, $false, $true, 1, $false, "83-79=", 2) | Out-Null
$d.Content.Find.Execute("51-45=", $true, $false, $false, $false, $false, $true, 1, $false, "77-30=", 2) | Out-Null
$d.Content.Find.Execute("11+45=", $true, $false, $false, $false, $false, $true, 1, $false, "11+25=", 2) | Out-Null
$d.Content.Find.Execute("17+25=", $true, $false, $false, $false, $false, $true, 1, $false, "30-26=", 2) | Out-Null
$d.Content.Find.Execute("23+9=", $true, $false, $false, $false, $false, $true, 1, $false, "85-11=", 2) | Out-Null
$d.Content.Find.Execute("13+30=", $true, $false, $false, $false, $false, $true, 1, $false, "70-58=", 2) | Out-Null
$d.Content.Find.Execute("65+32=", $true, $false, $false, $false, $false, $true, 1, $false, "86-65=", 2) | Out-Null
$d.Content.Find.Execute("32+1=", $true, $false, $false, $false, $false, $true, 1, $false, "12+74=", 2) | Out-Null
$d.Content.Find.Execute("33+11=", $true, $false, $false, $false, $false, $true, 1, $false, "11+46=", 2) | Out-Null
$d.Content.Find.Execute("68-28=", $true, $false, $false, $false, $false, $true, 1, $false, "24+66=", 2) | Out-Null
$d.Content.Find.Execute("29-3=", $true, $false, $false, $false, $false, $true, 1, $false, "86-71=", 2) | Out-Null
$d.Content.Find.Execute("32+46=", $true, $false, $false, $false, $false, $true, 1, $false, "89+7=", 2) | Out-Null
$d.Content.Find.Execute("10-6=", $true, $false, $false, $false, $false, $true, 1, $false, "95-80=", 2) | Out-Null
$d.Content.Find.Execute("27+43=", $true, $false, $false, $false, $false, $true, 1, $false, "13+58=", 2) | Out-Null
$d.Content.Find.Execute("1+35=", $true, $false, $false, $false, $false, $true, 1, $false, "75-44=", 2) | Out-Null
$d.Content.Find.Execute("34+32=", $true, $false, $false, $false, $false, $true, 1, $false, "76-73=", 2) | Out-Null
$d.Content.Find.Execute("99-9=", $true, $false, $false, $false, $false, $true, 1, $false, "89-38=", 2) | Out-Null
$d.Content.Find.Execute("5+54=", $true, $false, $false, $false, $false, $true, 1, $false, "71-63=", 2) | Out-Null
$d.Content.Find.Execute("95+4=", $true, $false, $false, $false, $false, $true, 1, $false, "47-11=", 2) | Out-Null
$d.Content.Find.Execute("58-34=", $true, $false, $false, $false, $false, $true, 1, $false, "91+7=", 2) | Out-Null
$d.Content.Find.Execute("7+59=", $true, $false, $false, $false, $false, $true, 1, $false, "4+53=", 2) | Out-Null
$d.Content.Find.Execute("33+64=", $true, $false, $false, $false, $false, $true, 1, $false, "59-6=", 2) | Out-Null
$d.Content.Find.Execute("98-57=", $true, $false, $false, $false, $false, $true, 1, $false, "27+68=", 2) | Out-Null
$d.Content.Find.Execute("4+38=", $true, $false, $false, $false, $false, $true, 1, $false, "56-34=", 2) | Out-Null
$d.Content.Find.Execute("1+64=", $true, $false, $false, $false, $false, $true, 1, $false, "17+80=", 2) | Out-Null
$d.Content.Find.Execute("55-11=", $true, $false, $false, $false, $false, $true, 1, $false, "42-25=", 2) | Out-Null
$d.Content.Find.Execute("16-13=", $true, $false, $false, $false, $false, $true, 1, $false, "85-14=", 2) | Out-Null
$d.Content.Find.Execute("35+59=", $true, $false, $false, $false, $false, $true, 1, $false, "78+0=", 2) | Out-Null
$d.Content.Find.Execute("20+17=", $true, $false, $false, $false, $false, $true, 1, $false, "64-33=", 2) | Out-Null
$d.Content.Find.Execute("1+61=", $true, $false, $false, $false, $false, $true, 1, $false, "49+17=", 2) | Out-Null
$d.Content.Find.Execute("10+48=", $true, $false, $false, $false, $false, $true, 1, $false, "46+34=", 2) | Out-Null
$d.Content.Find.Execute("23+38=", $true, $false, $false, $false, $false, $true, 1, $false, "50-20=", 2) | Out-Null
$d.Content.Find.Execute("89+5=", $true, $false, $false, $false, $false, $true, 1, $false, "68-59=", 2) | Out-Null
$d.Content.Find.Execute("56-4=", $true, $false, $false, $false, $false, $true, 1, $false, "5+68=", 2) | Out-Null
$d.Content.Find.Execute("19+47=", $true, $false, $false, $false, $false, $true, 1, $false, "75+14=", 2) | Out-Null
$d.Content.Find.Execute("53-48=", $true, $false, $false, $false, $false, $true, 1, $false, "67-12=", 2) | Out-Null
$d.Content.Find.Execute("77-74=", $true, $false, $false, $false, $false, $true, 1, $false, "33+6=", 2) | Out-Null
$d.Content.Find.Execute("70-41=", $true, $false, $false, $false, $false, $true, 1, $false, "99-67=", 2) | Out-Null
$d.Content.Find.Execute("37-19=", $true, $false, $false, $false, $false, $true, 1, $false, "33+12=", 2) | Out-Null
$d.Content.Find.Execute("40+58=", $true, $false, $false, $false, $false, $true, 1, $false, "49+30=", 2) | Out-Null
$d.Content.Find.Execute("90-45=", $true, $false, $false, $false, $false, $true, 1, $false, "1+47=", 2) | Out-Null
$d.Content.Find.Execute("92-88=", $true, $false, $false, $false, $false, $true, 1, $false, "91-51=", 2) | Out-Null
$d.Content.Find.Execute("42-40=", $true, $false, $false, $false, $false, $true, 1, $false, "8+42=", 2) | Out-Null
$d.Content.Find.Execute("35+55=", $true, $false, $false, $false, $false, $true, 1, $false, "27+63=", 2) | Out-Null
$d.Content.Find.Execute("17+1=", $true, $false, $false, $false, $false, $true, 1, $false, "1+17=", 2) | Out-Null
$d.Content.Find.Execute("70-69=", $true, $false, $false, $false, $false, $true, 1, $false, "98-44=", 2) | Out-Null
$d.Content.Find.Execute("79-50=", $true, $false, $false, $false, $false, $true, 1, $false, "50-34=", 2) | Out-Null
$d.Content.Find.Execute("14+54=", $true, $false, $false, $false, $false, $true, 1, $false, "59-2=", 2) | Out-Null
$d.Content.Find.Execute("92-13=", $true, $false, $false, $false, $false, $true, 1, $false, "70+20=", 2) | Out-Null
$d.Content.Find.Execute("0+54=", $true, $false, $false, $false, $false, $true, 1, $false, "66+30=", 2) | Out-Null
$d.Content.Find.Execute("40+55=", $true, $false, $false, $false, $false, $true, 1, $false, "64+18=", 2) | Out-Null
$d.Content.Find.Execute("51-8=", $true, $false, $false, $false, $false, $true, 1, $false, "74+21=", 2) | Out-Null
$d.Content.Find.Execute("96-21=", $true, $false, $false, $false, $false, $true, 1, $false, "84-34=", 2) | Out-Null
$d.Content.Find.Execute("44+39=", $true, $false, $false, $false, $false, $true, 1, $false, "44-12=", 2) | Out-Null
$d.Content.Find.Execute("57+1=", $true, $false, $false, $false, $false, $true, 1, $false, "20+52=", 2) | Out-Null
$d.Content.Find.Execute("73+14=", $true, $false, $false, $false, $false, $true, 1, $false, "17+28=", 2) | Out-Null
$d.Content.Find.Execute("30+59=", $true, $false, $false, $false, $false, $true, 1, $false, "83-61=", 2) | Out-Null
$d.Content.Find.Execute("99-16=", $true, $false, $false, $false, $false, $true, 1, $false, "53+13=", 2) | Out-Null
$d.Content.Find.Execute("87-1=", $true, $false, $false, $false, $false, $true, 1, $false, "54-27=", 2) | Out-Null
$d.Content.Find.Execute("47+13=", $true, $false, $false, $false, $false, $true, 1, $false, "64-32=", 2) | Out-Null
$d.Content.Find.Execute("9+54=", $true, $false, $false, $false, $false, $true, 1, $false, "58+20=", 2) | Out-Null
$d.Content.Find.Execute("19+72=", $true, $false, $false, $false, $false, $true, 1, $false, "15+47=", 2) | Out-Null
$d.Content.Find.Execute("51+5=", $true, $false, $false, $false, $false, $true, 1, $false, "30+40=", 2) | Out-Null
$d.Content.Find.Execute("67-29=", $true, $false, $false, $false, $false, $true, 1, $false, "74-23=", 2) | Out-Null
$d.Content.Find.Execute("40-33=", $true, $false, $false, $false, $false, $true, 1, $false, "42+53=", 2) | Out-Null
$d.Content.Find.Execute("56-32=", $true, $false, $false, $false, $false, $true, 1, $false, "1+37=", 2) | Out-Null
